$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 corresponds to the "DNN" model. Update the evaluation metrics
# to reflect the new values from the latest backup.
$ws.Range("B6").Value = 0.79
$ws.Range("C6").Value = 0.65
$ws.Range("E6").Value = 10
$ws.Range("F6").Value = 69
$ws.Range("H6").Value = 0.3
$ws.Range("J6").Value = 0.46
$ws.Range("K6").Value = 0.13
$ws.Range("L6").Value = 0.56
